$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.480.88'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '3.310.94'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("D5").Value = '''589.34'
$ws.Range("E5").Value = '  +2.81%  '
$ws.Range("D6").Value = '''180.32'
$ws.Range("E6").Value = '  +1.55%  '
$ws.Range("D7").Value = '''0.642'
$ws.Range("E7").Value = '  +2.05%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '3.314.93'
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").Value = '''0.126'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = '''6.86'
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("D12").Value = '''0.402'
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '3.887.24'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '66.507.58'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '''26.71'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '3.295.24'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").Value = '''426.83'
$ws.Range("E19").Value = '  -2.18%  '
$ws.Range("D20").Value = '''5.50'
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("D21").Value = '''13.10'
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("D22").Value = '''7.33'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '''71.38'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").Value = '''0.207'
$ws.Range("E27").Value = '  +6.43%  '
$ws.Range("D28").Value = '''0.0000115'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = '''9.46'
$ws.Range("E29").Value = '  +6.29%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").Value = '''22.39'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '''5.20'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").Value = '''6.61'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").Value = '''159.13'
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = '''1.44'
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("D39").Value = '2.864.88'
$ws.Range("E39").Value = '  +3.27%  '
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").Value = '''26.42'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").Value = '''4.36'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '''0.751'
$ws.Range("E43").Value = '  -3.68%  '
$ws.Range("D44").Value = '''39.74'
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = '''5.95'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").Value = '''2.32'
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("D47").Value = '''0.0639'
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").Value = '''315.16'
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = '''23.00'
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").Value = '''0.103'
$ws.Range("E51").Value = '  +0.51%  '

# Reset style on cells that required a quote-prefix to stay text,
# so no stray number-format / style is left behind on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
